# Avoid interactive "Are you sure you want to delete…" prompts when a
# pre-existing summary sheet has to be cleared out.
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Clear out any existing summary sheet(s) first so re-running this edit
# (or editing a workbook that already has one) doesn't leave stale /
# duplicate summary sheets behind.
foreach ($existing in @($wb.Worksheets)) {
    if ($existing.Name -eq "summary_F3") {
        $existing.Delete()
    }
}

# Add the new summary sheet after the last existing sheet (rawdata_F3),
# so it becomes the new, final, active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "summary_F3"

# Populate the single summary cell.
$newSheet.Range("A1").Value = "sdfds"

# Make sure the new summary sheet ends up the active/selected tab.
$newSheet.Activate()

$excel.DisplayAlerts = $true
